# Remove the "Device Stack / Device Drivers / ... / Hard Disk Drive
# Addresses" block of bullet paragraphs from the schedule, and strip the
# now-orphaned "Directory Data Structures" text (plus its trailing
# underline/bold line-break run) from the paragraph that carries the
# section break, leaving that paragraph empty but otherwise intact.

$d = $word.ActiveDocument

# Locate the paragraphs by their text, rather than a hard-coded index,
# so the edit is resilient to any incidental paragraph-count drift.
$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($startIdx -eq -1 -and $t -match "- Device Stack") {
        $startIdx = $i
    }
    if ($t -match "Directory Data Structures") {
        $endIdx = $i
    }
}

if ($startIdx -eq -1 -or $endIdx -eq -1) {
    throw "Could not locate target paragraphs (start=$startIdx end=$endIdx)"
}

$startPara = $d.Paragraphs.Item($startIdx)
$endPara = $d.Paragraphs.Item($endIdx)

# Delete every whole paragraph from "- Device Stack" up to (but not
# including) the "Directory Data Structures" paragraph. This removes:
#   - Device Stack
#   - Device Drivers
#     - Character Devices
#     - Block Devices
#     - Packet Devices
#   - Hard Disk Drive Addresses (C:H:S)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.Start)
$deleteRange.Delete()

# The "Directory Data Structures" paragraph (which holds the section
# break in its pPr) now sits where $endPara used to be. Clear its
# run content (the text plus the bold/underlined line-break run) while
# leaving the paragraph mark - and therefore the sectPr - in place.
$endPara = $d.Paragraphs.Item($startIdx)
$clearRange = $d.Range($endPara.Range.Start, $endPara.Range.End - 1)
if ($clearRange.Start -lt $clearRange.End) {
    $clearRange.Delete()
}
